# Convert numeric "count" cells into plain text cells (matching the
# formatting used by the rest of the sheet, which stores every value as
# text). Values >= 1000 get thousands separators (e.g. 2741 -> "2,741").
#
# Helper: write $Number into the cell at ($Row, $Col) of worksheet $Ws as
# literal text (no auto number/currency/percent reformatting by Excel).
function Set-TextNumber($Ws, $Row, $Col, $Number) {
    $text = "{0:N0}" -f $Number
    $cell = $Ws.Cells.Item($Row, $Col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Helper: write an arbitrary literal string into a cell as text, without
# Excel auto-detecting it as a currency/percentage/number.
function Set-TextValue($Ws, $Row, $Col, $Text) {
    $cell = $Ws.Cells.Item($Row, $Col)
    $cell.Value = "'" + $Text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overall": A2 2741 -> "2,741"
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextNumber $wsOverall 2 1 2741

# ---------------------------------------------------------------------
# Sheet "County": B2:B39 numeric -> text, plus new Total row 40
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")
$countyCounts = @(2,6,50,69,44,85,4,23,5,1,14,13,25,34,38,1129,91,21,9,23,2,15,29,13,5,209,36,51,4,168,185,13,102,1,33,104,13,72)
for ($i = 0; $i -lt $countyCounts.Length; $i++) {
    $row = $i + 2
    Set-TextNumber $wsCounty $row 2 $countyCounts[$i]
}

Set-TextValue $wsCounty 40 1 "Total"
Set-TextNumber $wsCounty 40 2 2741
Set-TextValue $wsCounty 40 3 "`$6,452,428,730"
Set-TextValue $wsCounty 40 4 "10.77%"
Set-TextValue $wsCounty 40 5 "-10.62%"
Set-TextValue $wsCounty 40 6 "66.00%"

# ---------------------------------------------------------------------
# Sheet "Congressional District": B2:B12 numeric -> text
# ---------------------------------------------------------------------
$wsCd = $wb.Worksheets.Item("Congressional District")
$cdCounts = @(149,160,336,152,193,263,356,618,166,348,2741)
for ($i = 0; $i -lt $cdCounts.Length; $i++) {
    $row = $i + 2
    Set-TextNumber $wsCd $row 2 $cdCounts[$i]
}

# ---------------------------------------------------------------------
# Sheet "Size": B2:B8 numeric -> text
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
$sizeCounts = @(849,810,514,167,303,98,2741)
for ($i = 0; $i -lt $sizeCounts.Length; $i++) {
    $row = $i + 2
    Set-TextNumber $wsSize $row 2 $sizeCounts[$i]
}

# ---------------------------------------------------------------------
# Sheet "Subsector": B2:B13 numeric -> text
# ---------------------------------------------------------------------
$wsSubsector = $wb.Worksheets.Item("Subsector")
$subsectorCounts = @(288,252,171,195,17,855,42,242,52,607,20,2741)
for ($i = 0; $i -lt $subsectorCounts.Length; $i++) {
    $row = $i + 2
    Set-TextNumber $wsSubsector $row 2 $subsectorCounts[$i]
}

Write-Output "Edit complete"
